# Remove the "INTERNAL" classification text-box/watermark shape that is
# anchored in every footer (primary, first-page and even-page) of the
# document. This corresponds to deleting the <w:r> run that wraps the
# mc:AlternateContent (drawing + VML fallback) inside each w:ftr part,
# leaving the empty "Fuzeile"-styled paragraph behind.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $footer = $sec.Footers.Item($i)
        if ($footer.Exists -or $footer.Shapes.Count -gt 0) {
            while ($footer.Shapes.Count -gt 0) {
                $footer.Shapes.Item(1).Delete()
            }
        }
    }
}
